$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report work: update Sunday hours for week of 43220 (row 16) from 0 to 1.
$ws.Range("C16").Value = 1

# UI tweaking: move the active selection.
$ws.Range("L22").Select()
